$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 165.28572
$ws.Range("I33").Value = 170.4
$ws.Range("K33").Value = 170.4
$ws.Range("M33").Value = 58.59999999999999
$ws.Range("H74").Value = 3311.4443
$ws.Range("I74").Value = 3100.5
$ws.Range("J74").Value = 3733.3333
$ws.Range("K74").Value = 3100.5
$ws.Range("L74").Value = 3733.3333
$ws.Range("M74").Value = -2164.5
$ws.Range("N74").Value = -5605.3333
$ws.Range("H77").Value = 3311.4443
$ws.Range("I77").Value = 3100.5
$ws.Range("J77").Value = 3733.3333
$ws.Range("K77").Value = 15502.5
$ws.Range("L77").Value = 18666.6665
$ws.Range("M77").Value = -10822.5
$ws.Range("N77").Value = -28026.6665
$ws.Range("H112").Value = 4959818
$ws.Range("J112").Value = 5455731.5
$ws.Range("L112").Value = 16367194.5
$ws.Range("N112").Value = -16369410.5
$ws.Range("H129").Value = 1112.3243
$ws.Range("I129").Value = 348.5
$ws.Range("J129").Value = 1204.909
$ws.Range("K129").Value = 1045.5
$ws.Range("L129").Value = 3614.727
$ws.Range("M129").Value = 3954.5
$ws.Range("N129").Value = -13614.727
$ws.Range("H132").Value = 150057.28
$ws.Range("I132").Value = 156431.94
$ws.Range("J132").Value = 25751.5
$ws.Range("K132").Value = 469295.82
$ws.Range("L132").Value = 77254.5
$ws.Range("M132").Value = -466765.82
$ws.Range("N132").Value = -82314.5
$ws.Range("H138").Value = 1427.68
$ws.Range("I138").Value = 765.74194
$ws.Range("J138").Value = 2507.6843
$ws.Range("K138").Value = 2297.22582
$ws.Range("L138").Value = 7523.0529
$ws.Range("M138").Value = 2842.77418
$ws.Range("N138").Value = -17803.0529
$ws.Range("H141").Value = 2349.397
$ws.Range("I141").Value = 1447.1724
$ws.Range("J141").Value = 7582.3
$ws.Range("K141").Value = 4341.5172
$ws.Range("L141").Value = 22746.9
$ws.Range("M141").Value = 838.4827999999998
$ws.Range("N141").Value = -33106.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2286.1292
$ws.Range("I61").Value = 1755.25
$ws.Range("J61").Value = 4106.2856
$ws.Range("K61").Value = 1755.25
$ws.Range("L61").Value = 4106.2856
$ws.Range("M61").Value = -1543.25
$ws.Range("N61").Value = -4530.2856
$ws.Range("H74").Value = 3899.46
$ws.Range("I74").Value = 1259.8158
$ws.Range("J74").Value = 12258.333
$ws.Range("K74").Value = 1259.8158
$ws.Range("L74").Value = 12258.333
$ws.Range("M74").Value = -385.8158000000001
$ws.Range("N74").Value = -14006.333
$ws.Range("H77").Value = 3899.46
$ws.Range("I77").Value = 1259.8158
$ws.Range("J77").Value = 12258.333
$ws.Range("K77").Value = 6299.079000000001
$ws.Range("L77").Value = 61291.665
$ws.Range("M77").Value = -1931.079000000001
$ws.Range("N77").Value = -70027.66500000001
$ws.Range("H132").Value = 2102.9492
$ws.Range("I132").Value = 2006.7693
$ws.Range("J132").Value = 2290.5
$ws.Range("K132").Value = 6020.3079
$ws.Range("L132").Value = 6871.5
$ws.Range("M132").Value = -3490.3079
$ws.Range("N132").Value = -11931.5
$ws.Range("H136").Value = 2286.1292
$ws.Range("I136").Value = 1755.25
$ws.Range("J136").Value = 4106.2856
$ws.Range("K136").Value = 5265.75
$ws.Range("L136").Value = 12318.8568
$ws.Range("M136").Value = -2715.75
$ws.Range("N136").Value = -17418.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6886.2383
$ws.Range("I86").Value = 2586.4614
$ws.Range("J86").Value = 13873.375
$ws.Range("K86").Value = 2586.4614
$ws.Range("L86").Value = 13873.375
$ws.Range("M86").Value = -1463.4614
$ws.Range("N86").Value = -16119.375
$ws.Range("H89").Value = 6886.2383
$ws.Range("I89").Value = 2586.4614
$ws.Range("J89").Value = 13873.375
$ws.Range("K89").Value = 12932.307
$ws.Range("L89").Value = 69366.875
$ws.Range("M89").Value = -7316.307000000001
$ws.Range("N89").Value = -80598.875
$ws.Range("H134").Value = 35718130
$ws.Range("I134").Value = 55558580
$ws.Range("J134").Value = 5311.6
$ws.Range("K134").Value = 166675740
$ws.Range("L134").Value = 15934.8
$ws.Range("M134").Value = -166673205
$ws.Range("N134").Value = -21004.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 29600
$ws.Range("J88").Value = 29600
$ws.Range("L88").Value = 29600
$ws.Range("N88").Value = -30412
$ws.Range("H91").Value = 29600
$ws.Range("J91").Value = 29600
$ws.Range("L91").Value = 29600
$ws.Range("N91").Value = -32408
$ws.Range("H132").Value = 1606.6875
$ws.Range("I132").Value = 1423.1163
$ws.Range("J132").Value = 3185.4
$ws.Range("K132").Value = 4269.3489
$ws.Range("L132").Value = 9556.2
$ws.Range("M132").Value = -1739.3489
$ws.Range("N132").Value = -14616.2
$ws.Range("H134").Value = 2010.9395
$ws.Range("I134").Value = 1365.4468
$ws.Range("K134").Value = 4096.3404
$ws.Range("M134").Value = -1561.3404

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1272.2325
$ws.Range("I5").Value = 920.3929
$ws.Range("J5").Value = 1929
$ws.Range("K5").Value = 2761.1787
$ws.Range("L5").Value = 5787
$ws.Range("M5").Value = -2649.1787
$ws.Range("N5").Value = -6011
$ws.Range("H131").Value = 1496.7441
$ws.Range("J131").Value = 1563.25
$ws.Range("L131").Value = 4689.75
$ws.Range("N131").Value = -14769.75
$ws.Range("H135").Value = 1272.2325
$ws.Range("I135").Value = 920.3929
$ws.Range("J135").Value = 1929
$ws.Range("K135").Value = 8283.536100000001
$ws.Range("L135").Value = 17361
$ws.Range("M135").Value = -5748.536100000001
$ws.Range("N135").Value = -22431

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3077.0444
$ws.Range("I132").Value = 2841.6177
$ws.Range("J132").Value = 3804.7273
$ws.Range("K132").Value = 8524.8531
$ws.Range("L132").Value = 11414.1819
$ws.Range("M132").Value = -5994.8531
$ws.Range("N132").Value = -16474.1819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6678.95
$ws.Range("I22").Value = 911
$ws.Range("J22").Value = 11398.182
$ws.Range("K22").Value = 911
$ws.Range("L22").Value = 11398.182
$ws.Range("M22").Value = -616
$ws.Range("N22").Value = -11988.182
$ws.Range("H27").Value = 6678.95
$ws.Range("I27").Value = 911
$ws.Range("J27").Value = 11398.182
$ws.Range("K27").Value = 911
$ws.Range("L27").Value = 11398.182
$ws.Range("M27").Value = -804
$ws.Range("N27").Value = -11612.182
$ws.Range("H132").Value = 6774.9536
$ws.Range("I132").Value = 6834.1934
$ws.Range("J132").Value = 6621.9165
$ws.Range("K132").Value = 20502.5802
$ws.Range("L132").Value = 19865.7495
$ws.Range("M132").Value = -17972.5802
$ws.Range("N132").Value = -24925.7495

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 10640622
$ws.Range("I132").Value = 15627111
$ws.Range("J132").Value = 2778.6667
$ws.Range("K132").Value = 46881333
$ws.Range("L132").Value = 8336.000100000001
$ws.Range("M132").Value = -46878803
$ws.Range("N132").Value = -13396.0001
